$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark from near the top of the document to its
#    new position inside the final section (between "mi" and "e abilità").
#    (We add it at the new spot further below; Bookmarks.Add with an existing
#    name moves the bookmark rather than duplicating it, so no separate
#    delete step is required.)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2. Rewrite the five bullet points in section 7 ("Riflessioni e
#    considerazioni conclusive"). The first bullet becomes a longer piece of
#    rich text (several runs plus the relocated _GoBack bookmark), and the
#    remaining bullets shift: each one now carries the text that used to
#    belong to the bullet above it. The last bullet's original text is
#    dropped entirely.
# ---------------------------------------------------------------------------

$bullet1 = $d.Paragraphs.Item(110).Range   # was: "Durante questo periodo di stage"
$bullet2 = $d.Paragraphs.Item(111).Range   # was: "indicare eventuali proposte/..."
$bullet3 = $d.Paragraphs.Item(112).Range   # was: "sono ipotizzabili ulteriori sviluppi/..."
$bullet4 = $d.Paragraphs.Item(113).Range   # was: "allegare, previa autorizzazione..."

# Cascade the plain-text bullets first (order doesn't matter since the new
# values are hard-coded strings, not derived from each other's live text).
$bullet4.Text = "sono ipotizzabili ulteriori sviluppi/varianti nel progetto/attività realizzate"
$bullet3.Text = "indicare eventuali proposte/suggerimenti per migliorare procedure o metodi di lavoro riscontrati in azienda"
$bullet2.Text = "Durante questo periodo di stage"

# Bullet 1 becomes the new rich-text sentence. Build it as plain text first
# (adjacent runs sharing identical formatting collapse into one <w:r>, which
# mirrors how Word itself coalesces runs), then splice the _GoBack bookmark
# in the middle of it.
$bullet1.Text = "Durante il periodo di stage ho potuto rafforzare le mie competenze comunicative, grazie anche alla collaborazione con figure facenti parte di altre aziende. Inoltre ho rafforzato le mie abilità organizzative e gestionali "

# Work out where "...ho rafforzato le mi" ends (right before "e abilità...")
# inside the now-updated bullet 1 paragraph, and drop the _GoBack bookmark
# there.
$p110 = $d.Paragraphs.Item(110).Range
$marker = "ho rafforzato le mi"
$bmPos = $p110.Start + $p110.Text.IndexOf($marker) + $marker.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3. Footer: the cached PAGE / NUMPAGES field results move from "2" and "6"
#    to "7" and "7" (the document grew by a page).
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)

$pageRng = $footer.Range.Duplicate
$pageRng.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "7", 1) | Out-Null

$numRng = $footer.Range.Duplicate
$numRng.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "7", 1) | Out-Null

Write-Output "edit applied"
